# Auto-generated edit script: updates recomputed market-price/profit
# columns (H:N) for the rows flagged by the scheduled pricing refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 18: You Grow, Girl / Growth Formula Beta
$ws.Range("H18").Value = 780
$ws.Range("I18").Value = 475
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 475
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -191
$ws.Range("N18").Value = -2568
# ALC row 32: Automata for the People / Crab Oil
$ws.Range("H32").Value = 2864.6667
$ws.Range("I32").Value = 4220.25
$ws.Range("J32").Value = 1780.2
$ws.Range("K32").Value = 4220.25
$ws.Range("L32").Value = 1780.2
$ws.Range("M32").Value = -3894.25
$ws.Range("N32").Value = -2432.2
# ALC row 101: Edge of the Arcane / Cunning Craftsman's Tea
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
# ALC row 127: Liquid Competence / Competent Craftsman's Draught
$ws.Range("H127").Value = 56433.707
$ws.Range("I127").Value = 56433.707
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 169301.121
$ws.Range("L127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -164341.121
# ALC row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 998.5
$ws.Range("I129").Value = 774.375
$ws.Range("J129").Value = 1895
$ws.Range("K129").Value = 2323.125
$ws.Range("L129").Value = 5685
$ws.Range("M129").Value = 2676.875
$ws.Range("N129").Value = -15685
# ALC row 131: Mindful Study / Grade 5 Tincture of Mind
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 48599.562
$ws.Range("I137").Value = 1915.5834
$ws.Range("J137").Value = 76609.95
$ws.Range("K137").Value = 5746.7502
$ws.Range("L137").Value = 229829.85
$ws.Range("M137").Value = -3196.7502
$ws.Range("N137").Value = -234929.85
# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4599.7
$ws.Range("I138").Value = 5554
$ws.Range("J138").Value = 4493.6665
$ws.Range("K138").Value = 16662
$ws.Range("L138").Value = 13480.9995
$ws.Range("M138").Value = -11522
$ws.Range("N138").Value = -23760.9995
# ALC row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 49218.473
$ws.Range("I141").Value = 49218.473
$ws.Range("K141").Value = 147655.419
$ws.Range("M141").Value = -142475.419

$ws = $wb.Worksheets.Item("ARM")
# ARM row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 199.71428
$ws.Range("I5").Value = 199.6
$ws.Range("K5").Value = 199.6
$ws.Range("M5").Value = -87.59999999999999
# ARM row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 139447.92
$ws.Range("I32").Value = 144905.14
$ws.Range("K32").Value = 144905.14
$ws.Range("M32").Value = -144618.14
# ARM row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 16605.092
$ws.Range("I45").Value = 17770.77
$ws.Range("J45").Value = 14921.333
$ws.Range("K45").Value = 17770.77
$ws.Range("L45").Value = 14921.333
$ws.Range("M45").Value = -17393.77
$ws.Range("N45").Value = -15675.333
# ARM row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2338.8
$ws.Range("I61").Value = 2237.9565
$ws.Range("J61").Value = 3498.5
$ws.Range("K61").Value = 2237.9565
$ws.Range("L61").Value = 3498.5
$ws.Range("M61").Value = -2025.9565
$ws.Range("N61").Value = -3922.5
# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 24043.166
$ws.Range("I132").Value = 25319.818
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 75959.454
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -73429.454
$ws.Range("N132").Value = -35060
# ARM row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2338.8
$ws.Range("I136").Value = 2237.9565
$ws.Range("J136").Value = 3498.5
$ws.Range("K136").Value = 6713.869499999999
$ws.Range("L136").Value = 10495.5
$ws.Range("M136").Value = -4163.869499999999
$ws.Range("N136").Value = -15595.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 199.71428
$ws.Range("I4").Value = 199.6
$ws.Range("K4").Value = 199.6
$ws.Range("M4").Value = -84.59999999999999
# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 8601.875
$ws.Range("I134").Value = 8137.8335
$ws.Range("J134").Value = 9994
$ws.Range("K134").Value = 24413.5005
$ws.Range("L134").Value = 29982
$ws.Range("M134").Value = -21878.5005
$ws.Range("N134").Value = -35052

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 453.75
$ws.Range("I22").Value = 288
$ws.Range("J22").Value = 951
$ws.Range("K22").Value = 288
$ws.Range("L22").Value = 951
$ws.Range("M22").Value = 62
$ws.Range("N22").Value = -1651
# CRP row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3261.5737
$ws.Range("I31").Value = 2626.28
$ws.Range("J31").Value = 3702.75
$ws.Range("K31").Value = 2626.28
$ws.Range("L31").Value = 3702.75
$ws.Range("M31").Value = -2331.28
$ws.Range("N31").Value = -4292.75
# CRP row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3261.5737
$ws.Range("I34").Value = 2626.28
$ws.Range("J34").Value = 3702.75
$ws.Range("K34").Value = 2626.28
$ws.Range("L34").Value = 3702.75
$ws.Range("M34").Value = -2424.28
$ws.Range("N34").Value = -4106.75

$ws = $wb.Worksheets.Item("CUL")
# CUL row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 10006000
$ws.Range("J68").Value = 12506250
$ws.Range("L68").Value = 37518750
$ws.Range("N68").Value = -37520372
# CUL row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 10006000
$ws.Range("J71").Value = 12506250
$ws.Range("L71").Value = 112556250
$ws.Range("N71").Value = -112564362
# CUL row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value = 127981.19
$ws.Range("I129").Value = 334423
$ws.Range("J129").Value = 4116.1
$ws.Range("K129").Value = 1003269
$ws.Range("L129").Value = 12348.3
$ws.Range("M129").Value = -998269
$ws.Range("N129").Value = -22348.3
# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 128569.34
$ws.Range("I131").Value = 60027.65
$ws.Range("J131").Value = 258037
$ws.Range("K131").Value = 180082.95
$ws.Range("L131").Value = 774111
$ws.Range("M131").Value = -175042.95
$ws.Range("N131").Value = -784191
# CUL row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 4548879
$ws.Range("J137").Value = 4086.75
$ws.Range("L137").Value = 12260.25
$ws.Range("N137").Value = -22460.25

$ws = $wb.Worksheets.Item("GSM")
# GSM row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 26734.16
$ws.Range("I97").Value = 47352.645
$ws.Range("J97").Value = 492.45456
$ws.Range("K97").Value = 47352.645
$ws.Range("L97").Value = 492.45456
$ws.Range("M97").Value = -46856.645
$ws.Range("N97").Value = -1484.45456

$ws = $wb.Worksheets.Item("LTW")
# LTW row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 18836.525
$ws.Range("I93").Value = 1041.7
$ws.Range("J93").Value = 38608.555
$ws.Range("K93").Value = 1041.7
$ws.Range("L93").Value = 38608.555
$ws.Range("M93").Value = 206.3
$ws.Range("N93").Value = -41104.555

$ws = $wb.Worksheets.Item("WVR")
# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 7613.4062
$ws.Range("I132").Value = 8797.083000000001
$ws.Range("K132").Value = 26391.249
$ws.Range("M132").Value = -23861.249
# WVR row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 7405.8184
$ws.Range("I136").Value = 7646.4
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 22939.2
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -20389.2
$ws.Range("N136").Value = -20100
